# "Primeros cambios con Yendry"
# Repurpose the "Retiro" (withdrawal) receipt template as a
# "Cancelación" (cancellation) receipt template:
#   - Title cell: "Comprobante de Retiro" -> "Comprobante de Cancelación"
#   - Sub-title cell: "Retiro" -> "Cancelación"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Comprobante de Cancelación"
$ws.Range("A23").Value = "Cancelación"

# Leave the selection on the merged "Cancelación" cell, matching where
# the user's edit session ended up.
$ws.Range("A23:D23").Select()
